$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 32
$ws1.Range("H2").Value = 0.31
$ws1.Range("L2").Value = 0.82

$ws1.Range("D3").Value = 33
$ws1.Range("L3").Value = 1.06

$ws1.Range("D4").Value = 32
$ws1.Range("L4").Value = 0.97

$ws1.Range("L5").Value = 1.02

$ws1.Range("D6").Value = 30
$ws1.Range("L6").Value = 0.99

$ws1.Range("D7").Value = 29
$ws1.Range("L7").Value = 0.9399999999999999

$ws1.Range("L8").Value = 0.85

$ws1.Range("L9").Value = 1.1

$ws1.Range("D10").Value = 27
$ws1.Range("L10").Value = 0.9

$ws1.Range("D11").Value = 27
$ws1.Range("L11").Value = 1.09

$ws1.Range("D12").Value = 27
$ws1.Range("L12").Value = 0.91

$ws1.Range("L13").Value = 0.9399999999999999

$ws1.Range("D14").Value = 27
$ws1.Range("L14").Value = 0.92

$ws1.Range("D15").Value = 27
$ws1.Range("L15").Value = 0.88

$ws1.Range("D16").Value = 28
$ws1.Range("L16").Value = 1.08

$ws1.Range("D17").Value = 27
$ws1.Range("L17").Value = 1.03

# --- Sheet "Summary" ---
# B9:B14 hold numeric-looking values stored as TEXT in the workbook
# (e.g. "422", "213", ...). Force text format before assigning so the
# runtime doesn't silently coerce these into numeric cells.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "459"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "242"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "128"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "33"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "27"
